# TrialsSetup.xlsx update (2026-02-10 12:00)
# - The OPERA-2 trial row was removed from the query results (trial
#   completed/closed out), shifting all rows below it up by one.
# - REJOICE (MK-5909-003) days remaining dropped from 10 to 9.
# - REMASTER (CLOU) days remaining dropped from 30 to 29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# OPERA-2 is the second row of data (sheet row 2) - remove it entirely so
# the table/list object and every row below it shift up.
$ws.Rows(2).Delete()

# After the shift, REJOICE is now row 6 and REMASTER is now row 8.
$ws.Range("B6").Value = 9
$ws.Range("B8").Value = 29

# Keep the workbook-level ExternalData_1 defined name (used by the query
# table) in sync with the table's new, smaller extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!ExternalData_1") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$9"
    }
}
